$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-weight values (rows 2-17) per author revision
$rowData = @{
    2 = @{ "E"=3; "G"=50.63667066666667; "H"=151.910012; "I"=0.20081482031288; "J"=0.20081482031288; "K"=3; "M"=5.015232333333333; "N"=15.045697; "O"=0.2232512241152976; "P"=0.2232512241152976; "Q"=253.9546679798182; "R"=2285.592011818364; "S"=0.044832154455344; "T"=0.044832154455344 }
    3 = @{ "E"=3; "G"=50.63667066666667; "H"=151.910012; "I"=0.20081482031288; "J"=0.20081482031288; "K"=3; "M"=7.971374; "N"=23.914122; "O"=0.3548427839629211; "P"=0.3548427839629211; "Q"=403.6438399988293; "R"=3632.794559989464; "S"=0.07125768990083611; "T"=0.0712576899008361 }
    4 = @{ "E"=3; "G"=50.63667066666667; "H"=151.910012; "I"=0.20081482031288; "J"=0.20081482031288; "K"=3; "M"=6.516197000000001; "N"=19.548591; "O"=0.290066114615979; "P"=0.290066114615979; "Q"=329.9585214881214; "R"=2969.626693393092; "S"=0.05824957468546308; "T"=0.05824957468546307 }
    5 = @{ "E"=3; "G"=50.63667066666667; "H"=151.910012; "I"=0.20081482031288; "J"=0.20081482031288; "K"=3; "M"=2.96172; "N"=8.885159999999999; "O"=0.1318398773058023; "P"=0.1318398773058023; "Q"=149.97164024688; "R"=1349.74476222192; "S"=0.02647540127123684; "T"=0.02647540127123684 }
    6 = @{ "E"=3; "G"=67.324; "H"=201.972; "I"=0.2669934019110801; "J"=0.2669934019110801; "K"=3; "M"=5.015232333333333; "N"=15.045697; "O"=0.2232512241152976; "P"=0.2232512241152976; "Q"=337.6455016093333; "R"=3038.809514484; "S"=0.05960660380735628; "T"=0.05960660380735628 }
    7 = @{ "E"=3; "G"=67.324; "H"=201.972; "I"=0.2669934019110801; "J"=0.2669934019110801; "K"=3; "M"=7.971374; "N"=23.914122; "O"=0.3548427839629211; "P"=0.3548427839629211; "Q"=536.664783176; "R"=4829.983048584; "S"=0.09474068203385877; "T"=0.09474068203385876 }
    8 = @{ "E"=3; "G"=67.324; "H"=201.972; "I"=0.2669934019110801; "J"=0.2669934019110801; "K"=3; "M"=6.516197000000001; "N"=19.548591; "O"=0.290066114615979; "P"=0.290066114615979; "Q"=438.696446828; "R"=3948.268021452001; "S"=0.0774457387204495; "T"=0.0774457387204495 }
    9 = @{ "E"=3; "G"=67.324; "H"=201.972; "I"=0.2669934019110801; "J"=0.2669934019110801; "K"=3; "M"=2.96172; "N"=8.885159999999999; "O"=0.1318398773058023; "P"=0.1318398773058023; "Q"=199.39483728; "R"=1794.55353552; "S"=0.03520037734941557; "T"=0.03520037734941557 }
    10 = @{ "E"=3; "G"=71.72398199999999; "H"=215.171946; "I"=0.2844428428612245; "J"=0.2844428428612245; "K"=3; "M"=5.015232333333333; "N"=15.045697; "O"=0.2232512241152976; "P"=0.2232512241152976; "Q"=359.7124336018179; "R"=3237.411902416361; "S"=0.06350221285960361; "T"=0.06350221285960361 }
    11 = @{ "E"=3; "G"=71.72398199999999; "H"=215.171946; "I"=0.2844428428612245; "J"=0.2844428428612245; "K"=3; "M"=7.971374; "N"=23.914122; "O"=0.3548427839629211; "P"=0.3548427839629211; "Q"=571.738685291268; "R"=5145.648167621412; "S"=0.1009324902392046; "T"=0.1009324902392046 }
    12 = @{ "E"=3; "G"=71.72398199999999; "H"=215.171946; "I"=0.2844428428612245; "J"=0.2844428428612245; "K"=3; "M"=6.516197000000001; "N"=19.548591; "O"=0.290066114615979; "P"=0.290066114615979; "Q"=467.367596336454; "R"=4206.308367028087; "S"=0.08250723025907883; "T"=0.08250723025907883 }
    13 = @{ "E"=3; "G"=71.72398199999999; "H"=215.171946; "I"=0.2844428428612245; "J"=0.2844428428612245; "K"=3; "M"=2.96172; "N"=8.885159999999999; "O"=0.1318398773058023; "P"=0.1318398773058023; "Q"=212.42635196904; "R"=1911.83716772136; "S"=0.03750090950333744; "T"=0.03750090950333744 }
    14 = @{ "E"=3; "G"=62.47139133333334; "H"=187.414174; "I"=0.2477489349148154; "J"=0.2477489349148154; "K"=3; "M"=5.015232333333333; "N"=15.045697; "O"=0.2232512241152976; "P"=0.2232512241152976; "Q"=313.3085417232531; "R"=2819.776875509278; "S"=0.05531025299299375; "T"=0.05531025299299375 }
    15 = @{ "E"=3; "G"=62.47139133333334; "H"=187.414174; "I"=0.2477489349148154; "J"=0.2477489349148154; "K"=3; "M"=7.971374; "N"=23.914122; "O"=0.3548427839629211; "P"=0.3548427839629211; "Q"=497.9828246183587; "R"=4481.845421565228; "S"=0.08791192178902167; "T"=0.08791192178902164 }
    16 = @{ "E"=3; "G"=62.47139133333334; "H"=187.414174; "I"=0.2477489349148154; "J"=0.2477489349148154; "K"=3; "M"=6.516197000000001; "N"=19.548591; "O"=0.290066114615979; "P"=0.290066114615979; "Q"=407.0758927920928; "R"=3663.683035128834; "S"=0.07186357095098757; "T"=0.07186357095098757 }
    17 = @{ "E"=3; "G"=62.47139133333334; "H"=187.414174; "I"=0.2477489349148154; "J"=0.2477489349148154; "K"=3; "M"=2.96172; "N"=8.885159999999999; "O"=0.1318398773058023; "P"=0.1318398773058023; "Q"=185.02276913976; "R"=1665.20492225784; "S"=0.03266318918181248; "T"=0.03266318918181247 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}